$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# B-column values that differ from 0 for the newly appended rows (59-79)
$bValues = @{ 61 = 24; 78 = 22 }

$startRow = 59
$endRow = 79
$startSerial = 44273

for ($r = $startRow; $r -le $endRow; $r++) {
    $serial = $startSerial + ($r - $startRow)
    $ws.Cells.Item($r, 1).Value = $serial

    $bVal = 0
    if ($bValues.ContainsKey($r)) { $bVal = $bValues[$r] }
    $ws.Cells.Item($r, 2).Value = $bVal

    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
}

# Match the date number format used by the existing date column (A2:A58)
$ws.Cells.Item($startRow - 1, 1).Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view to reflect the new selection/scroll position
try {
    $excel.ActiveWindow.ScrollRow = 52
} catch {
}
$ws.Range("B61").Select()
